$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
Write-Host $ws.Name
